$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.415.06'
$ws.Range("E2").Value = '  -1.18%  '

# Row 3
$ws.Range("D3").Value = '1.711.53'
$ws.Range("E3").Value = '  -1.49%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.11%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '224.46'
$ws.Range("E5").Value = '  -1.47%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5340'
$ws.Range("E6").Value = '  -2.17%  '

# Row 7
$ws.Range("E7").Value = '  +0.14%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2672'
$ws.Range("E8").Value = '  -3.15%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06615'

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.93'
$ws.Range("E10").Value = '  -4.44%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07643'
$ws.Range("E11").Value = '  -1.86%  '

# Row 12
$ws.Range("E12").Value = '  -2.93%  '

# Row 13
$ws.Range("D13").Value = '1.709.22'
$ws.Range("E13").Value = '  -2.64%  '

# Row 14
$ws.Range("D14").Value = '1.948.28'
$ws.Range("E14").Value = '  -1.44%  '

# Row 15
$ws.Range("E15").Value = '  -3.71%  '

# Row 16
$ws.Range("D16").Value = '0.0₅8163'
$ws.Range("E16").Value = '  -3.15%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.90'
$ws.Range("E17").Value = '  -2.17%  '

# Row 18
$ws.Range("D18").Value = '27.393.56'
$ws.Range("E18").Value = '  -1.28%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '216.54'

# Row 20
$ws.Range("E20").Value = '  +0.11%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.663'
$ws.Range("E21").Value = '  -3.76%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.46'

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.977'
$ws.Range("E23").Value = '  -4.18%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.006'
$ws.Range("E24").Value = '  +0.09%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '142.57'
$ws.Range("E25").Value = '  -3.19%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.731'
$ws.Range("E26").Value = '  +0.06%  '

# Row 27
$ws.Range("E27").Value = '  -2.84%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.286'
$ws.Range("E28").Value = '  -2.35%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '16.26'
$ws.Range("E29").Value = '  -5.40%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05401'
$ws.Range("E30").Value = '  -4.86%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.291'
$ws.Range("E31").Value = '  -1.78%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.492'
$ws.Range("E32").Value = '  -5.68%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.428'
$ws.Range("E33").Value = '  -2.56%  '

# Row 34
$ws.Range("E34").Value = '  -2.64%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.881'
$ws.Range("E35").Value = '  +0.89%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9488'

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.415'
$ws.Range("E37").Value = '  -1.40%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5840'
$ws.Range("E38").Value = '  -2.29%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01633'
$ws.Range("E39").Value = '  -2.18%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.862'
$ws.Range("E40").Value = '  -1.04%  '

# Row 41
$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.005'
$ws.Range("E41").Value = '  +0.15%  '

# Row 42
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '1.044.16'
$ws.Range("E42").Value = '  -0.64%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8394'
$ws.Range("E43").Value = '  -1.21%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.85'
$ws.Range("E44").Value = '  -1.21%  '

# Row 45
$ws.Range("D45").Value = '1.854.75'
$ws.Range("E45").Value = '  -1.40%  '

# Row 46
$ws.Range("E46").Value = '  +1.94%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '57.97'
$ws.Range("E47").Value = '  -2.76%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4518'
$ws.Range("E48").Value = '  +1.78%  '

# Row 49
$ws.Range("E49").Value = '  +0.13%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.047'
$ws.Range("E50").Value = '  -3.28%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05243'
$ws.Range("E51").Value = '  -1.60%  '
